# Fixed bug in beacon creater
# Adds an example bit-row (row 12) and the packed beacon value (row 14)
# to the worksheet, updating the used range/selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: example bit values for columns C..T (18 bits)
$row12Values = @(1,1,0,1,1,0,0,1,0,0,0,0,0,1,0,0,0,1)
$col = 3  # column C
foreach ($v in $row12Values) {
    $ws.Cells.Item(12, $col).Value = $v
    $col++
}

# Row 14: packed numeric value of the beacon
$ws.Range("C14").Value = 1820467200

# Update selection to reflect the new "next empty" row, as Excel would.
$ws.Range("C15").Select()
